$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = 0.8673805046080227
$ws.Range("F1").Value = -1.570796280873402

$ws.Range("E2").Value = 0.8693143401741139
$ws.Range("F2").Value = -1.570796286790826

$ws.Range("E3").Value = 0.87797920561113
$ws.Range("F3").Value = -1.570796313304808

$ws.Range("E4").Value = 0.8901698485995828
$ws.Range("F4").Value = -1.570796350607462

$ws.Range("E5").Value = 0.8988347140365989
$ws.Range("F5").Value = -1.570796377121443

$ws.Range("E6").Value = 0.90076854960269
$ws.Range("F6").Value = -1.570796383038867
